$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(30)
$r = $p.Range
$r.Collapse(0)

# Insert four empty paragraphs after "Try to make it engaging..." -- each
# InsertParagraphAfter() splits off a new "Normal"-style paragraph, matching
# the target markup (<w:p><w:r><w:t>...</w:t></w:r></w:p>, no explicit pPr).
$r.InsertParagraphAfter()
$r2 = $d.Paragraphs.Item(31).Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r3 = $d.Paragraphs.Item(32).Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$r4 = $d.Paragraphs.Item(33).Range
$r4.Collapse(0)
$r4.InsertParagraphAfter()

$d.Paragraphs.Item(31).Range.Text = "Making sure your code doesn’t have issues"
$d.Paragraphs.Item(32).Range.Text = "Now you want to make sure you code doesn’t have any errors in it.  First, let’s make sure it compiles without any warnings.  In Visual Studio, choose Build>Rebuild Solution from the menu and make sure the error list at the bottom of the window doesn’t have any errors or compiler warnings."
$d.Paragraphs.Item(33).Range.Text = "Now go to the Unity window and find the Unity “Console”.  You’ll find it in the Console tab in the bottom pane of the window.  This is where exceptions get printed if your code throws and exception.  You can also display messages here manually using Unity’s Debug.Log() method.  However, the final code you turn in for your project should not call Debug.Log() or otherwise print any messages in the console window."
$d.Paragraphs.Item(34).Range.Text = "Run your project.  Let it run for a minute or so, pressing buttons and moving the joysticks around, just to make sure no errors happen and you don’t have any Debug.Log() calls left."

# Only the first of the four (the new section heading) should be Heading1;
# set it last so the style doesn't propagate onto paragraphs split off later.
$d.Paragraphs.Item(31).Style = "Heading1"

for ($i = 29; $i -le 36; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Output "$i [$($pp.Style.NameLocal)]: $($pp.Range.Text)"
}
